# Weekly update: a new week's record is inserted at the top of the data
# (row 6), pushing the existing historical rows down by one and adding a
# new row 24 at the bottom (previously row 23's data).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 6; this shifts rows 6:23 down to 7:24 and
# extends the sheet dimension to A1:R24 automatically, carrying styles
# (e.g. the date format on column D) along with the shifted rows.
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with the new week's data.
$ws.Cells.Item(6, 1).Value2 = 7
$ws.Cells.Item(6, 2).Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(6, 3).Value2 = "Ñuble"
$ws.Cells.Item(6, 4).Value2 = 44980
$ws.Cells.Item(6, 5).Value2 = 16
$ws.Cells.Item(6, 6).Value2 = 100112012
$ws.Cells.Item(6, 7).Value2 = "Espinaca"
$ws.Cells.Item(6, 8).Value2 = "Sin especificar"
$ws.Cells.Item(6, 9).Value2 = "Primera"
$ws.Cells.Item(6, 10).Value2 = 60
$ws.Cells.Item(6, 11).Value2 = 7500
$ws.Cells.Item(6, 12).Value2 = 8000
$ws.Cells.Item(6, 13).Value2 = 7750
$ws.Cells.Item(6, 14).Value2 = "$/cuna 10 kilos"
$ws.Cells.Item(6, 15).Value2 = "Provincia de Diguillín"
$ws.Cells.Item(6, 16).Value2 = 775
$ws.Cells.Item(6, 17).Value2 = 10
$ws.Cells.Item(6, 18).Value2 = "Hortaliza"
